$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume figures (column D = Price, column E =
# Volume(1h) change). Values that look numeric are written with the cell
# pre-set to Text format so Excel stores them as literal strings (matching
# the source data feed, which keeps thousand-separator-less decimal
# formatting such as trailing zeros, e.g. "246.40") instead of coercing
# them into numeric cells.
$updates = @(
    @{ Cell = 'D2'; Value = '30.560.53'; Text = $false },
    @{ Cell = 'E2'; Value = '  -0.58%  '; Text = $false },
    @{ Cell = 'D3'; Value = '1.884.57'; Text = $false },
    @{ Cell = 'E4'; Value = '  -0.04%  '; Text = $false },
    @{ Cell = 'D5'; Value = '246.40'; Text = $true },
    @{ Cell = 'E5'; Value = '  -0.85%  '; Text = $false },
    @{ Cell = 'E6'; Value = '  -0.06%  '; Text = $false },
    @{ Cell = 'D7'; Value = '0.4741'; Text = $true },
    @{ Cell = 'E7'; Value = '  +0.03%  '; Text = $false },
    @{ Cell = 'D8'; Value = '0.2894'; Text = $true },
    @{ Cell = 'E8'; Value = '  -1.06%  '; Text = $false },
    @{ Cell = 'D9'; Value = '0.06545'; Text = $true },
    @{ Cell = 'E9'; Value = '  +0.15%  '; Text = $false },
    @{ Cell = 'D10'; Value = '22.35'; Text = $true },
    @{ Cell = 'E10'; Value = '  +1.05%  '; Text = $false },
    @{ Cell = 'D11'; Value = '0.7763'; Text = $true },
    @{ Cell = 'E11'; Value = '  +5.21%  '; Text = $false },
    @{ Cell = 'D12'; Value = '101.09'; Text = $true },
    @{ Cell = 'E12'; Value = '  +4.35%  '; Text = $false },
    @{ Cell = 'D13'; Value = '0.07815'; Text = $true },
    @{ Cell = 'E13'; Value = '  +0.14%  '; Text = $false },
    @{ Cell = 'D14'; Value = '1.883.04'; Text = $false },
    @{ Cell = 'E14'; Value = '  -0.30%  '; Text = $false },
    @{ Cell = 'D15'; Value = '5.259'; Text = $true },
    @{ Cell = 'E15'; Value = '  +0.29%  '; Text = $false },
    @{ Cell = 'D16'; Value = '284.91'; Text = $true },
    @{ Cell = 'E16'; Value = '  +0.01%  '; Text = $false },
    @{ Cell = 'D17'; Value = '30.546.81'; Text = $false },
    @{ Cell = 'E17'; Value = '  -0.57%  '; Text = $false },
    @{ Cell = 'D18'; Value = '13.23'; Text = $true },
    @{ Cell = 'E18'; Value = '  -0.27%  '; Text = $false },
    @{ Cell = 'D19'; Value = '0.000007536'; Text = $true },
    @{ Cell = 'E19'; Value = '  -0.17%  '; Text = $false },
    @{ Cell = 'E20'; Value = '  +0.01%  '; Text = $false },
    @{ Cell = 'D21'; Value = '2.129.66'; Text = $false },
    @{ Cell = 'E21'; Value = '  -0.28%  '; Text = $false },
    @{ Cell = 'D22'; Value = '5.363'; Text = $true },
    @{ Cell = 'E22'; Value = '  +0.80%  '; Text = $false },
    @{ Cell = 'E23'; Value = '  -0.08%  '; Text = $false },
    @{ Cell = 'D24'; Value = '6.468'; Text = $true },
    @{ Cell = 'E24'; Value = '  +3.42%  '; Text = $false },
    @{ Cell = 'D25'; Value = '9.172'; Text = $true },
    @{ Cell = 'E25'; Value = '  -0.59%  '; Text = $false },
    @{ Cell = 'D26'; Value = '163.20'; Text = $true },
    @{ Cell = 'E26'; Value = '  -1.10%  '; Text = $false },
    @{ Cell = 'D27'; Value = '19.14'; Text = $true },
    @{ Cell = 'E27'; Value = '  +0.76%  '; Text = $false },
    @{ Cell = 'D28'; Value = '1.916'; Text = $true },
    @{ Cell = 'E28'; Value = '  -0.17%  '; Text = $false },
    @{ Cell = 'D29'; Value = '1.335'; Text = $true },
    @{ Cell = 'D30'; Value = '0.09706'; Text = $true },
    @{ Cell = 'E30'; Value = '  -0.36%  '; Text = $false },
    @{ Cell = 'D31'; Value = '1.503'; Text = $true },
    @{ Cell = 'E31'; Value = '  +0.41%  '; Text = $false },
    @{ Cell = 'D32'; Value = '4.251'; Text = $true },
    @{ Cell = 'E32'; Value = '  -1.10%  '; Text = $false },
    @{ Cell = 'D33'; Value = '4.193'; Text = $true },
    @{ Cell = 'E33'; Value = '  +0.14%  '; Text = $false },
    @{ Cell = 'D34'; Value = '0.04851'; Text = $true },
    @{ Cell = 'E34'; Value = '  -0.22%  '; Text = $false },
    @{ Cell = 'D35'; Value = '1.131'; Text = $true },
    @{ Cell = 'E35'; Value = '  +0.43%  '; Text = $false },
    @{ Cell = 'D36'; Value = '0.6990'; Text = $true },
    @{ Cell = 'E36'; Value = '  +0.23%  '; Text = $false },
    @{ Cell = 'D37'; Value = '2.761'; Text = $true },
    @{ Cell = 'E37'; Value = '  +1.36%  '; Text = $false },
    @{ Cell = 'D38'; Value = '0.01918'; Text = $true },
    @{ Cell = 'E38'; Value = '  +1.30%  '; Text = $false },
    @{ Cell = 'D39'; Value = '2.901'; Text = $true },
    @{ Cell = 'E39'; Value = '  +3.34%  '; Text = $false },
    @{ Cell = 'D40'; Value = '76.05'; Text = $true },
    @{ Cell = 'E40'; Value = '  -0.08%  '; Text = $false },
    @{ Cell = 'D41'; Value = '6.299'; Text = $true },
    @{ Cell = 'E41'; Value = '  -0.47%  '; Text = $false },
    @{ Cell = 'D42'; Value = '1.992'; Text = $true },
    @{ Cell = 'E42'; Value = '  -0.14%  '; Text = $false },
    @{ Cell = 'D43'; Value = '0.4258'; Text = $true },
    @{ Cell = 'E43'; Value = '  -0.42%  '; Text = $false },
    @{ Cell = 'E44'; Value = '  -0.13%  '; Text = $false },
    @{ Cell = 'D45'; Value = '0.8317'; Text = $true },
    @{ Cell = 'E45'; Value = '  -0.43%  '; Text = $false },
    @{ Cell = 'D46'; Value = '101.61'; Text = $true },
    @{ Cell = 'E46'; Value = '  -0.10%  '; Text = $false },
    @{ Cell = 'D47'; Value = '9.889'; Text = $true },
    @{ Cell = 'E47'; Value = '  +3.35%  '; Text = $false },
    @{ Cell = 'D48'; Value = '7.029'; Text = $true },
    @{ Cell = 'E48'; Value = '  -0.27%  '; Text = $false },
    @{ Cell = 'D49'; Value = '35.22'; Text = $true },
    @{ Cell = 'E49'; Value = '  -1.13%  '; Text = $false },
    @{ Cell = 'D50'; Value = '896.16'; Text = $true },
    @{ Cell = 'E50'; Value = '  -2.05%  '; Text = $false },
    @{ Cell = 'D51'; Value = '0.05770'; Text = $true },
    @{ Cell = 'E51'; Value = '  +0.21%  '; Text = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Text) {
        $cell.NumberFormat = '@'
    }
    $cell.Value = $u.Value
}
